$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 12500
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 23000
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = -1532
$ws.Range("N21").Value = -23936

# Row 23
$ws.Range("H23").Value = 12500
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 23000
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 23000
$ws.Range("M23").Value = -1766
$ws.Range("N23").Value = -23468

# Row 34
$ws.Range("H34").Value = 1310.7273
$ws.Range("I34").Value = 1310.7273
$ws.Range("K34").Value = 1310.7273
$ws.Range("M34").Value = -1107.7273

# Row 36
$ws.Range("H36").Value = 1310.7273
$ws.Range("I36").Value = 1310.7273
$ws.Range("K36").Value = 1310.7273
$ws.Range("M36").Value = -595.7273

# Row 86
$ws.Range("H86").Value = 55558864
$ws.Range("I86").Value = 166667500
$ws.Range("J86").Value = 4546.6665
$ws.Range("K86").Value = 166667500
$ws.Range("L86").Value = 4546.6665
$ws.Range("M86").Value = -166666377
$ws.Range("N86").Value = -6792.6665

# Row 89
$ws.Range("H89").Value = 55558864
$ws.Range("I89").Value = 166667500
$ws.Range("J89").Value = 4546.6665
$ws.Range("K89").Value = 833337500
$ws.Range("L89").Value = 22733.3325
$ws.Range("M89").Value = -833331884
$ws.Range("N89").Value = -33965.3325

# Row 116
$ws.Range("H116").Value = 1874.375
$ws.Range("I116").Value = 1882.5
$ws.Range("J116").Value = 1850
$ws.Range("K116").Value = 1882.5
$ws.Range("L116").Value = 1850
$ws.Range("M116").Value = 1559.5
$ws.Range("N116").Value = -8734

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8449.52
$ws.Range("I32").Value = 7492.9507
$ws.Range("J32").Value = 12527.526
$ws.Range("K32").Value = 7492.9507
$ws.Range("L32").Value = 12527.526
$ws.Range("M32").Value = -7205.9507
$ws.Range("N32").Value = -13101.526

# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# Row 61
$ws.Range("H61").Value = 15627028
$ws.Range("I61").Value = 17859056
$ws.Range("K61").Value = 17859056
$ws.Range("M61").Value = -17858844

# Row 110
$ws.Range("H110").Value = 478.22223
$ws.Range("I110").Value = 529.1429000000001
$ws.Range("J110").Value = 300
$ws.Range("K110").Value = 529.1429000000001
$ws.Range("L110").Value = 300
$ws.Range("M110").Value = 1515.8571
$ws.Range("N110").Value = -4390

# Row 136
$ws.Range("H136").Value = 15627028
$ws.Range("I136").Value = 17859056
$ws.Range("K136").Value = 53577168
$ws.Range("M136").Value = -53574618

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2047.05
$ws.Range("I20").Value = 2103.8
$ws.Range("J20").Value = 1876.8
$ws.Range("K20").Value = 2103.8
$ws.Range("L20").Value = 1876.8
$ws.Range("M20").Value = -1856.8
$ws.Range("N20").Value = -2370.8

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 94
$ws.Range("H94").Value = 688.29034
$ws.Range("I94").Value = 639.8570999999999
$ws.Range("J94").Value = 790
$ws.Range("K94").Value = 639.8570999999999
$ws.Range("L94").Value = 790
$ws.Range("M94").Value = -188.8570999999999
$ws.Range("N94").Value = -1692

# Row 99
$ws.Range("H99").Value = 1333.3334
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 1250
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1250
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -4246

# Row 107
$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 755.55554
$ws.Range("K107").Value = 755.55554
$ws.Range("M107").Value = 1164.44446

# Row 134
$ws.Range("H134").Value = 3629.853
$ws.Range("I134").Value = 2961.7307
$ws.Range("J134").Value = 5801.25
$ws.Range("K134").Value = 8885.1921
$ws.Range("L134").Value = 17403.75
$ws.Range("M134").Value = -6350.1921
$ws.Range("N134").Value = -22473.75

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 628465.4
$ws.Range("I134").Value = 1843.4482
$ws.Range("J134").Value = 2647580.5
$ws.Range("K134").Value = 5530.3446
$ws.Range("L134").Value = 7942741.5
$ws.Range("M134").Value = -2995.3446
$ws.Range("N134").Value = -7947811.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 992144.3
$ws.Range("J2").Value = 2314823
$ws.Range("L2").Value = 13888938
$ws.Range("N2").Value = -13889164

# Row 117
$ws.Range("H117").Value = 488
$ws.Range("I117").Value = 216.66667
$ws.Range("J117").Value = 895
$ws.Range("K117").Value = 650.00001
$ws.Range("L117").Value = 2685
$ws.Range("M117").Value = 2791.99999
$ws.Range("N117").Value = -9569

$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 70
$ws.Range("H70").Value = 15406.294
$ws.Range("I70").Value = 52625
$ws.Range("J70").Value = 3954.3845
$ws.Range("K70").Value = 52625
$ws.Range("L70").Value = 3954.3845
$ws.Range("M70").Value = -52355
$ws.Range("N70").Value = -4494.3845

# Row 73
$ws.Range("H73").Value = 15406.294
$ws.Range("I73").Value = 52625
$ws.Range("J73").Value = 3954.3845
$ws.Range("K73").Value = 52625
$ws.Range("L73").Value = 3954.3845
$ws.Range("M73").Value = -51689
$ws.Range("N73").Value = -5826.3845

# Row 132
$ws.Range("H132").Value = 4941.8184
$ws.Range("I132").Value = 5538.871
$ws.Range("J132").Value = 3518.077
$ws.Range("K132").Value = 16616.613
$ws.Range("L132").Value = 10554.231
$ws.Range("M132").Value = -14086.613
$ws.Range("N132").Value = -15614.231

$ws = $wb.Worksheets.Item("LTW")
# Row 45
$ws.Range("H45").Value = 5500
$ws.Range("I45").Value = 5500
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5500
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -5093
$ws.Range("N45").ClearContents()

# Row 93
$ws.Range("H93").Value = 1500
$ws.Range("I93").Value = 1500
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1500
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -252
$ws.Range("N93").ClearContents()

# Row 132
$ws.Range("H132").Value = 10208036
$ws.Range("I132").Value = 3565.1072
$ws.Range("J132").Value = 23813998
$ws.Range("K132").Value = 10695.3216
$ws.Range("L132").Value = 71441994
$ws.Range("M132").Value = -8165.321599999999
$ws.Range("N132").Value = -71447054

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 4755.5
$ws.Range("I20").Value = 4250
$ws.Range("K20").Value = 4250
$ws.Range("M20").Value = -4010

# Row 132
$ws.Range("H132").Value = 1282.2858
$ws.Range("I132").Value = 1099.6666
$ws.Range("J132").Value = 1680.7273
$ws.Range("K132").Value = 3298.9998
$ws.Range("L132").Value = 5042.1819
$ws.Range("M132").Value = -768.9998000000001
$ws.Range("N132").Value = -10102.1819
